$d = $word.ActiveDocument

# The edit removes everything in the document body after the first
# paragraph ("Centrix Web-Scraping Application"), i.e. all the scraped
# web content (text + inline PNG images) that followed it, leaving only
# that initial heading paragraph before the section properties.

if ($d.Paragraphs.Count -gt 1) {
    $start = $d.Paragraphs(2).Range.Start
    $end = $d.Content.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
